# "shoulder & base 2.0" -- rework of the CAD parameters sheet.
# The parameter list is renamed/reordered/re-valued wholesale, so we just
# rewrite every row's label (A), value (B) and unit (C) to match the new
# design rather than patching individual cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = "base_length"
$ws.Cells.Item(1, 2).Value = 7
$ws.Cells.Item(1, 3).Value = "in"

$ws.Cells.Item(2, 1).Value = "base_width"
$ws.Cells.Item(2, 2).Value = 5
$ws.Cells.Item(2, 3).Value = "in"

$ws.Cells.Item(3, 1).Value = "base_height"
$ws.Cells.Item(3, 2).Value = 0.5
$ws.Cells.Item(3, 3).Value = "in"

$ws.Cells.Item(4, 1).Value = "base_cutout_contact_angle"
$ws.Cells.Item(4, 2).Value = 60
$ws.Cells.Item(4, 3).Value = "deg"

$ws.Cells.Item(5, 1).Value = "base_cutout_diameter"
$ws.Cells.Item(5, 2).Value = 1.75
$ws.Cells.Item(5, 3).Value = "in"

$ws.Cells.Item(6, 1).Value = "turntable_motor_support_thickness"
$ws.Cells.Item(6, 2).Value = 0.125
$ws.Cells.Item(6, 3).Value = "in"

$ws.Cells.Item(7, 1).Value = "turntable_motor_support_width"
$ws.Cells.Item(7, 2).Value = 1.875
$ws.Cells.Item(7, 3).Value = "in"

$ws.Cells.Item(8, 1).Value = "turntable_motor_support_height"
$ws.Cells.Item(8, 2).Value = 0.75
$ws.Cells.Item(8, 3).Value = "in"

$ws.Cells.Item(9, 1).Value = "turntable_motor_boss"
$ws.Cells.Item(9, 2).Value = 1
$ws.Cells.Item(9, 3).Value = "in"

$ws.Cells.Item(10, 1).Value = "turntable_bed_depth"
$ws.Cells.Item(10, 2).Value = 0.125
$ws.Cells.Item(10, 3).Value = "in"

$ws.Cells.Item(11, 1).Value = "turntable_diameter"
$ws.Cells.Item(11, 2).Value = 4
$ws.Cells.Item(11, 3).Value = "in"

$ws.Cells.Item(12, 1).Value = "turntable_height"
$ws.Cells.Item(12, 2).Value = 0.375
$ws.Cells.Item(12, 3).Value = "in"

$ws.Cells.Item(13, 1).Value = "turntable_spurs"
$ws.Cells.Item(13, 2).Value = 6
$ws.Cells.Item(13, 3).Value = "ul"

$ws.Cells.Item(14, 1).Value = "turntable_spur_width"
$ws.Cells.Item(14, 2).Value = 0.25
$ws.Cells.Item(14, 3).Value = "in"

$ws.Cells.Item(15, 1).Value = "turntable_retention_shaft_diameter"
$ws.Cells.Item(15, 2).Value = 0.25
$ws.Cells.Item(15, 3).Value = "in"

$ws.Cells.Item(16, 1).Value = "spine_diameter"
$ws.Cells.Item(16, 2).Value = 1.25
$ws.Cells.Item(16, 3).Value = "in"

$ws.Cells.Item(17, 1).Value = "spine_height"
$ws.Cells.Item(17, 2).Value = 1
$ws.Cells.Item(17, 3).Value = "in"

$ws.Cells.Item(18, 1).Value = "shoulder_length"
$ws.Cells.Item(18, 2).Value = 2
$ws.Cells.Item(18, 3).Value = "in"

$ws.Cells.Item(19, 1).Value = "shoulder_width"
$ws.Cells.Item(19, 2).Value = 2.75
$ws.Cells.Item(19, 3).Value = "in"

$ws.Cells.Item(20, 1).Value = "shoulder_wall_thickness"
$ws.Cells.Item(20, 2).Value = 0.1875
$ws.Cells.Item(20, 3).Value = "in"

$ws.Cells.Item(21, 1).Value = "shoulder_gear_diameter"
$ws.Cells.Item(21, 2).Value = 2.5
$ws.Cells.Item(21, 3).Value = "in"

$ws.Cells.Item(22, 1).Value = "shoulder_motor_gear_diameter"
$ws.Cells.Item(22, 2).Value = 1
$ws.Cells.Item(22, 3).Value = "in"

$ws.Cells.Item(23, 1).Value = "shoulder_gear_shaft_diameter"
$ws.Cells.Item(23, 2).Value = 0.5
$ws.Cells.Item(23, 3).Value = "in"

$ws.Cells.Item(24, 1).Value = "shoulder_motor_boss"
$ws.Cells.Item(24, 2).Value = 1
$ws.Cells.Item(24, 3).Value = "in"

# Column A grew wider to fit the longer parameter names.
$ws.Columns.Item(1).ColumnWidth = 40.7109375

# The author's cursor ended up parked on the spur-count cell.
$ws.Range("B13").Select() | Out-Null
